# "Improved documentation and test data"
#
# The "Data" sheet is a denormalized user-story table (Role | Goal | Reason |
# Acceptance Criteria). This edit:
#   1. Nudges the saved window position back to the left edge of the screen.
#   2. Inserts a new first (header/example) row with a little-joke sample
#      story ("As a Beatle, I want to hold your hand, so that I have a hit
#      single. Acceptance criteria: the song must be less than three minutes
#      long") used to illustrate the table format.
#   3. Shortens the wording of the "quickly and easily test my code..." goal
#      used by the very first real data row (now row 2), dropping the
#      redundant trailing clause.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: move the saved window back to x=0 (matches xWindow 7700 -> 0).
$wb.Windows.Item(1).Left = 0
$excel.ActiveWindow.Left = 0

# Push all the existing rows down by one and make room for the new example
# row at the very top of the table.
$ws.Rows.Item(1).Insert()

# New illustrative example row (row 1).
$ws.Range("A1").Value = "Beatle"
$ws.Range("B1").Value = "hold your hand"
$ws.Range("D1").Value = "The song must be less than three minutes long"
$ws.Range("C1").Value = "have a hit single"

# Trim the wording on what is now row 2 (previously row 1).
$ws.Range("B2").Value = "quickly and easily test my code in my test Cloud environment"

# Put the selection/cursor on the first real data row's goal cell.
[void]$ws.Range("B2").Select()
